$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1): two new date columns (Jun_17, Jun_15) are inserted
# before the existing dates, pushing the old Jun_13 (B1) and Jun_10 (C1)
# out to D1 and E1.
$ws.Range("E1").Value = $ws.Range("C1").Value()
$ws.Range("D1").Value = $ws.Range("B1").Value()
$ws.Range("C1").Value = "Jun_15"
$ws.Range("B1").Value = "Jun_17"

# --- Data rows (2-27): two new "UN" rating columns are appended after the
# existing B/C columns (A/B/C stay exactly where they were).
$ws.Range("D2:E27").Value = "UN"

# --- Column widths: new columns D and E pick up the same fixed width as
# the existing custom-width column C (8 characters stored width).
$ws.Columns.Item(3).ColumnWidth = 7.166666666666667
$ws.Columns.Item(4).ColumnWidth = 7.166666666666667
$ws.Columns.Item(5).ColumnWidth = 7.166666666666667
